# Update the cryptos list worksheet with refreshed price / volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep the literal text representation (avoids Excel
    # coercing dotted-number-looking strings into floating point numbers),
    # then restore the default "Normal" style so no stray number format is
    # left behind on the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Price column (D) updates ---
Set-TextValue "D2"  "26.643.93"
Set-TextValue "D3"  "1.642.36"
Set-TextValue "D5"  "215.35"
Set-TextValue "D9"  "0.0627"
Set-TextValue "D10" "19.28"
Set-TextValue "D12" "1.871.85"

# Row 13 / Row 14 swap (coin order changed, with new price & volume figures)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.646.49"
$ws.Range("E13").Value = "  +0.80%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "4.19"
$ws.Range("E14").Value = "  +1.93%  "

Set-TextValue "D16" "65.44"
Set-TextValue "D17" "26.697.27"
Set-TextValue "D18" "0.0₃0744"
Set-TextValue "D23" "9.51"
Set-TextValue "D25" "145.44"
Set-TextValue "D28" "7.18"
Set-TextValue "D29" "15.76"
Set-TextValue "D32" "3.38"
Set-TextValue "D34" "1.278.73"
Set-TextValue "D37" "2.41"
Set-TextValue "D38" "0.533"
Set-TextValue "D41" "0.815"
Set-TextValue "D44" "1.782.24"
Set-TextValue "D45" "92.08"
Set-TextValue "D46" "59.86"
Set-TextValue "D49" "7.81"

# --- Volume(1h) column (E) updates ---
$ws.Range("E3").Value  = "  +0.46%  "
$ws.Range("E4").Value  = "  +0.11%  "
$ws.Range("E5").Value  = "  +0.75%  "
$ws.Range("E6").Value  = "  +1.26%  "
$ws.Range("E7").Value  = "  +0.14%  "
$ws.Range("E8").Value  = "  +0.03%  "
$ws.Range("E9").Value  = "  +0.71%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("E24").Value = "  +13.77%  "
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("E36").Value = "  +4.74%  "
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("E38").Value = "  +6.33%  "
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("E46").Value = "  +7.94%  "
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("E51").Value = "  -0.49%  "
